# "Basic Mouseover for Inventory"
# Update the Advanced-Inventory task row (row 18) and the row below it
# (row 19) with new effort/remaining numbers, then leave the selection
# on C19 (matching the author's final cursor position).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18: "Advanced Inventory (e.g.: Pop-Up on Moseover, Item Selection, ...)"
#   Curr. Est. (h): 6 -> 9
#   Effort (h):     1.5 -> 4
# (Remain(h) and Completion(%) are formulas and recalc automatically.)
$ws.Range("C18").Value = 9
$ws.Range("D18").Value = 4

# Row 19: "Aufgabe auf Sascha übertragen."
#   Curr. Est. (h): 2 -> 3
#   Effort (h):     2 -> 3
$ws.Range("C19").Value = 3
$ws.Range("D19").Value = 3

# Match the saved selection/active cell from the authored workbook.
$ws.Range("C19").Select()
